$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5589519839034186
$ws.Range("D2").Value = 0.5818401712027534

$ws.Range("C3").Value = 1.10156351950153
$ws.Range("D3").Value = 0.2825577250161191

$ws.Range("C4").Value = 0.1244755502189791
$ws.Range("D4").Value = 0.9020691155412068

$ws.Range("C5").Value = 0.1765663356847198
$ws.Range("D5").Value = 0.8614643740339849

$ws.Range("C6").Value = 0.3716140265417159
$ws.Range("D6").Value = 0.7137353326550464

$ws.Range("C7").Value = -0.8086275764389266
$ws.Range("D7").Value = 0.4273878994668798

$ws.Range("C8").Value = -0.473365998010179
$ws.Range("D8").Value = 0.6406154285286039

$ws.Range("C9").Value = -1.093104675759802
$ws.Range("D9").Value = 0.2861732832637562

$ws.Range("C10").Value = -0.6366810953078901
$ws.Range("D10").Value = 0.5309021439155921

$ws.Range("C11").Value = 0.09514687176493136
$ws.Range("D11").Value = 0.9250595036899281

$wb.Save()
